# Generate Report for Handoff
# b.md has now been handed off again (new xliff files generated) for both
# the zh-cn and de-de locales. Update the Overview sheet and each locale
# sheet's row for b.md to reflect the new "Ready for handoff" status, the
# new handoff xliff file names/timestamps, and the version-mismatch error
# that was detected against the (now out of date) handback file.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e329d65c1d2dd70e658bbf803e07ab06af65dd96/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/30100ab8eb8441d515c4d2882e29c4641ad6c290/e2e/b.md."

# ---- Overview sheet: row 3 is b.md ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-07 06:54:39"

# ---- zh-cn sheet: row 3 is b.md ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe keeps this a literal text "False" (matching the other
# text-typed True/False cells in the sheet) instead of Excel's automatic
# TRUE/FALSE boolean-literal coercion.
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-07 06:54:34"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.1667

# ---- de-de sheet: row 3 is b.md ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-07 06:54:39"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.1667
